# Generate Report for Handback
# Refresh the handback-status report with the latest handoff/handback
# timestamps and priority values produced by a subsequent report run.
$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# "Latest HO Xliff Generate Date" for the 3adca88e.md / 4539e755.md rows
# (rows 2 and 3 share the same text, so both must be updated together)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-04 18:19:57"
$wsOverview.Range("G3").Value = "2016-09-04 18:19:57"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-09-04 18:19:52"
$wsZhCn.Range("H3").Value = "2016-09-04 18:19:52"
$wsZhCn.Range("K2").Value = "2016-09-04 18:20:18"
$wsZhCn.Range("K3").Value = "2016-09-04 18:20:18"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
# H (Correspond Handoff Datetime) shares its text with Overview!G2:G3
$wsDeDe.Range("H2").Value = "2016-09-04 18:19:57"
$wsDeDe.Range("H3").Value = "2016-09-04 18:19:57"
# K (Correspond Handback DateTime)
$wsDeDe.Range("K2").Value = "2016-09-04 18:20:26"
$wsDeDe.Range("K3").Value = "2016-09-04 18:20:26"
